# Updated cryptos list with latest price/volume data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextValue "D2" '69.328.42'
Set-TextValue "E2" '  +2.14%  '
Set-TextValue "D3" '3.393.81'
Set-TextValue "E3" '  +1.60%  '
Set-TextValue "D4" '1.00'
Set-TextValue "D5" '587.66'
Set-TextValue "E5" '  +0.73%  '
Set-TextValue "D6" '180.82'
Set-TextValue "E6" '  +2.20%  '
Set-TextValue "E7" '  +0.03%  '
Set-TextValue "E8" '  +0.72%  '
Set-TextValue "E9" '  +8.21%  '
Set-TextValue "D10" '0.590'
Set-TextValue "E10" '  +1.62%  '
Set-TextValue "D11" '48.90'
Set-TextValue "E11" '  +4.58%  '
Set-TextValue "D12" '0.0000284'
Set-TextValue "E12" '  +4.14%  '
Set-TextValue "D13" '683.20'
Set-TextValue "E13" '  -1.18%  '
Set-TextValue "D14" '8.65'
Set-TextValue "E14" '  +2.21%  '
Set-TextValue "D15" '3.940.75'
Set-TextValue "E15" '  +1.46%  '
Set-TextValue "D16" '69.418.14'
Set-TextValue "E16" '  +2.22%  '
Set-TextValue "D17" '3.398.87'
Set-TextValue "E17" '  +1.84%  '
Set-TextValue "E18" '  +1.73%  '
Set-TextValue "D19" '17.72'
Set-TextValue "E19" '  +1.71%  '
Set-TextValue "D20" '11.39'
Set-TextValue "E20" '  +2.56%  '
Set-TextValue "D21" '0.902'
Set-TextValue "E21" '  +0.75%  '
Set-TextValue "E22" '  +1.21%  '
Set-TextValue "D23" '17.13'
Set-TextValue "E23" '  +0.49%  '
Set-TextValue "D24" '103.89'
Set-TextValue "E24" '  +5.32%  '
Set-TextValue "E25" '  +1.08%  '
Set-TextValue "D26" '2.74'
Set-TextValue "E26" '  +1.36%  '
Set-TextValue "D27" '9.63'
Set-TextValue "E27" '  +1.15%  '
Set-TextValue "D28" '34.26'
Set-TextValue "E28" '  +3.93%  '
Set-TextValue "D29" '8.74'
Set-TextValue "D30" '6.99'
Set-TextValue "E30" '  -1.89%  '
$ws.Range("B31").Value = 'Cosmos'
$ws.Range("C31").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextValue "D31" '11.21'
Set-TextValue "E31" '  +1.82%  '
$ws.Range("B32").Value = 'dogwifhat'
$ws.Range("C32").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
Set-TextValue "D32" '3.69'
Set-TextValue "E32" '  +9.58%  '
Set-TextValue "D33" '555.88'
Set-TextValue "E33" '  -2.94%  '
Set-TextValue "E34" '  +0.77%  '
Set-TextValue "E35" '  +1.44%  '
Set-TextValue "E36" '  +0.01%  '
Set-TextValue "D37" '3.704.19'
Set-TextValue "E37" '  -0.49%  '
Set-TextValue "D38" '0.140'
Set-TextValue "E38" '  +7.00%  '
Set-TextValue "D39" '35.10'
Set-TextValue "E39" '  +2.75%  '
Set-TextValue "D40" '3.25'
Set-TextValue "E40" '  +1.64%  '
Set-TextValue "D41" '0.0₃0707'
Set-TextValue "E41" '  +4.73%  '
Set-TextValue "E42" '  +1.16%  '
Set-TextValue "D43" '0.340'
Set-TextValue "E43" '  +1.06%  '
Set-TextValue "D44" '0.0424'
Set-TextValue "E44" '  +4.20%  '
Set-TextValue "E45" '  -2.42%  '
Set-TextValue "E46" '  -0.19%  '
Set-TextValue "E47" '  +0.77%  '
$ws.Range("B48").Value = 'Mantle'
$ws.Range("C48").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
Set-TextValue "D48" '1.39'
Set-TextValue "E48" '  +4.78%  '
$ws.Range("B49").Value = 'FirstDigitalUSD'
$ws.Range("C49").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
Set-TextValue "D49" '1.00'
Set-TextValue "E49" '  +0.01%  '
Set-TextValue "D50" '132.02'
Set-TextValue "E50" '  +1.72%  '
$ws.Range("B51").Value = 'CoreDAO'
$ws.Range("C51").Value = 'https://coinranking.com/coin/HFvoXUQh4+coredao-core'
Set-TextValue "D51" '2.58'
Set-TextValue "E51" '  -1.92%  '
